$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A98").Value = "2025-04-29 16:51:32"
$ws.Range("B98").Value = 276
